$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''310.78'
$ws.Range("E2").Value = '''-0.79%'
$ws.Range("D3").Value = '''36.81'
$ws.Range("E3").Value = '''-2.28%'
$ws.Range("D4").Value = '''5.110'
$ws.Range("E4").Value = '''-0.29%'
$ws.Range("D5").Value = '''0.07774'
$ws.Range("E5").Value = '''-1.82%'
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '''8.334'
$ws.Range("E6").Value = '''0.68%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.861'
$ws.Range("E7").Value = '''-3.76%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.948'
$ws.Range("E8").Value = '''1.01%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9297'
$ws.Range("E9").Value = '''1.04%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1151'
$ws.Range("E10").Value = '''-6.48%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1891'
$ws.Range("E11").Value = '''-2.02%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.08829'
$ws.Range("E12").Value = '''-4.45%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03281'
$ws.Range("E13").Value = '''-0.93%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09479'
$ws.Range("E14").Value = '''-1.38%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001372'
$ws.Range("E15").Value = '''-1.08%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005812'
$ws.Range("E16").Value = '''0.41%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.348'
$ws.Range("E17").Value = '''-4.51%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''4.370'
$ws.Range("E18").Value = '''-1.00%'
$ws.Range("D19").Value = '''0.3435'
$ws.Range("E19").Value = '''-0.27%'
$ws.Range("D20").Value = '''6.286'
$ws.Range("E20").Value = '''19.28%'
$ws.Range("D21").Value = '''0.1288'
$ws.Range("E21").Value = '''1.20%'
$ws.Range("E22").Value = '''-10.74%'
$ws.Range("D23").Value = '''0.04350'
$ws.Range("E23").Value = '''-0.09%'
$ws.Range("D24").Value = '''0.001195'
$ws.Range("E24").Value = '''-4.30%'
$ws.Range("D25").Value = '''0.004291'
$ws.Range("E25").Value = '''-0.42%'
$ws.Range("D26").Value = '''0.0001326'
$ws.Range("E26").Value = '''8.82%'
$ws.Range("D27").Value = '''0.0002880'
$ws.Range("E27").Value = '''-98.63%'
$ws.Range("D39").Value = '''0.02145'
$ws.Range("E39").Value = '''-4.19%'
$ws.Range("D40").Value = '''0.05069'
$ws.Range("E40").Value = '''-1.01%'
$ws.Range("D41").Value = '''0.007495'
$ws.Range("E41").Value = '''0.51%'
$ws.Range("D42").Value = '''0.1346'
$ws.Range("E42").Value = '''-1.32%'
$ws.Range("D43").Value = '''0.008394'
$ws.Range("E43").Value = '''-4.45%'
$ws.Range("D44").Value = '''0.001990'
$ws.Range("E44").Value = '''-0.91%'
$ws.Range("D45").Value = '''0.007796'
$ws.Range("E45").Value = '''-9.33%'
$ws.Range("D46").Value = '''0.00006317'
$ws.Range("E46").Value = '''-6.18%'
$ws.Range("D47").Value = '''0.00000000745'
$ws.Range("E47").Value = '''-0.60%'
$ws.Range("D48").Value = '''0.002859'
$ws.Range("E48").Value = '''-14.59%'
$ws.Range("D49").Value = '''0.001678'
$ws.Range("E49").Value = '''39.90%'
$ws.Range("D50").Value = '''0.00002086'
$ws.Range("E50").Value = '''-0.60%'
$ws.Range("D51").Value = '''0.0001986'
$ws.Range("E51").Value = '''-0.60%'
